# chore: update Sheets via scheduled runner
# Refreshes cached market-board price / profit figures (columns H-N) on a
# handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 50002668
$ws.Range("I40").Value = 2896.4
$ws.Range("J40").Value = 100002440
$ws.Range("K40").Value = 2896.4
$ws.Range("L40").Value = 100002440
$ws.Range("M40").Value = -2721.4
$ws.Range("N40").Value = -100002790

$ws.Range("H43").Value = 3535.9048
$ws.Range("I43").Value = 3717.3076
$ws.Range("J43").Value = 3241.125
$ws.Range("K43").Value = 3717.3076
$ws.Range("L43").Value = 3241.125
$ws.Range("M43").Value = -3648.3076
$ws.Range("N43").Value = -3379.125

$ws.Range("H98").Value = 2342.5833
$ws.Range("J98").Value = 4067.6667
$ws.Range("L98").Value = 4067.6667
$ws.Range("N98").Value = -7063.6667

$ws.Range("H122").Value = 2342.5833
$ws.Range("J122").Value = 4067.6667
$ws.Range("L122").Value = 12203.0001
$ws.Range("N122").Value = -17103.0001

$ws.Range("H132").Value = 3341.4688
$ws.Range("I132").Value = 2497.5667
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 7492.7001
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -4962.7001
$ws.Range("N132").Value = -53060

$ws.Range("H137").Value = 2476.44
$ws.Range("I137").Value = 2396.1365
$ws.Range("K137").Value = 7188.4095
$ws.Range("M137").Value = -4638.4095

$ws.Range("H138").Value = 4420.41
$ws.Range("I138").Value = 4421.625
$ws.Range("J138").Value = 4420.1787
$ws.Range("K138").Value = 13264.875
$ws.Range("L138").Value = 13260.5361
$ws.Range("M138").Value = -8124.875
$ws.Range("N138").Value = -23540.5361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 8999.5
$ws.Range("I28").Value = 8999.5
$ws.Range("K28").Value = 8999.5
$ws.Range("M28").Value = -8807.5

$ws.Range("H32").Value = 21075.633
$ws.Range("I32").Value = 20433.576
$ws.Range("J32").Value = 25249
$ws.Range("K32").Value = 20433.576
$ws.Range("L32").Value = 25249
$ws.Range("M32").Value = -20146.576
$ws.Range("N32").Value = -25823

$ws.Range("H61").Value = 12361894
$ws.Range("I61").Value = 28586500
$ws.Range("J61").Value = 1004670.8
$ws.Range("K61").Value = 28586500
$ws.Range("L61").Value = 1004670.8
$ws.Range("M61").Value = -28586288
$ws.Range("N61").Value = -1005094.8

$ws.Range("H74").Value = 1958.3
$ws.Range("I74").Value = 2121.111
$ws.Range("J74").Value = 493
$ws.Range("K74").Value = 2121.111
$ws.Range("L74").Value = 493
$ws.Range("M74").Value = -1247.111
$ws.Range("N74").Value = -2241

$ws.Range("H77").Value = 1958.3
$ws.Range("I77").Value = 2121.111
$ws.Range("J77").Value = 493
$ws.Range("K77").Value = 10605.555
$ws.Range("L77").Value = 2465
$ws.Range("M77").Value = -6237.555
$ws.Range("N77").Value = -11201

$ws.Range("H88").Value = 3262.2942
$ws.Range("I88").Value = 2399.6667
$ws.Range("J88").Value = 3732.818
$ws.Range("K88").Value = 2399.6667
$ws.Range("L88").Value = 3732.818
$ws.Range("M88").Value = -1993.6667
$ws.Range("N88").Value = -4544.818

$ws.Range("H91").Value = 3262.2942
$ws.Range("I91").Value = 2399.6667
$ws.Range("J91").Value = 3732.818
$ws.Range("K91").Value = 2399.6667
$ws.Range("L91").Value = 3732.818
$ws.Range("M91").Value = -995.6667000000002
$ws.Range("N91").Value = -6540.818

$ws.Range("H99").Value = 8999.5
$ws.Range("I99").Value = 8999.5
$ws.Range("K99").Value = 8999.5
$ws.Range("M99").Value = -6004.5

$ws.Range("H136").Value = 12361894
$ws.Range("I136").Value = 28586500
$ws.Range("J136").Value = 1004670.8
$ws.Range("K136").Value = 85759500
$ws.Range("L136").Value = 3014012.4
$ws.Range("M136").Value = -85756950
$ws.Range("N136").Value = -3019112.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7825.636
$ws.Range("J86").Value = 8335.625
$ws.Range("L86").Value = 8335.625
$ws.Range("N86").Value = -10581.625

$ws.Range("H89").Value = 7825.636
$ws.Range("J89").Value = 8335.625
$ws.Range("L89").Value = 41678.125
$ws.Range("N89").Value = -52910.125

$ws.Range("H94").Value = 2217
$ws.Range("I94").Value = 2435.55
$ws.Range("K94").Value = 2435.55
$ws.Range("M94").Value = -1984.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7693319
$ws.Range("I16").Value = 9091668
$ws.Range("K16").Value = 9091668
$ws.Range("M16").Value = -9091381

$ws.Range("H28").Value = 82965.336
$ws.Range("J28").Value = 82965.336
$ws.Range("L28").Value = 82965.336
$ws.Range("N28").Value = -83455.336

$ws.Range("H31").Value = 23259436
$ws.Range("I31").Value = 25644146
$ws.Range("J31").Value = 8524.75
$ws.Range("K31").Value = 25644146
$ws.Range("L31").Value = 8524.75
$ws.Range("M31").Value = -25643851
$ws.Range("N31").Value = -9114.75

$ws.Range("H34").Value = 23259436
$ws.Range("I34").Value = 25644146
$ws.Range("J34").Value = 8524.75
$ws.Range("K34").Value = 25644146
$ws.Range("L34").Value = 8524.75
$ws.Range("M34").Value = -25643944
$ws.Range("N34").Value = -8928.75

$ws.Range("H37").Value = 9020.5
$ws.Range("J37").Value = 11990
$ws.Range("L37").Value = 11990
$ws.Range("N37").Value = -12204

$ws.Range("H103").Value = 26552.1
$ws.Range("J103").Value = 56317.332
$ws.Range("L103").Value = 56317.332
$ws.Range("N103").Value = -58661.332

$ws.Range("H113").Value = 7693319
$ws.Range("I113").Value = 9091668
$ws.Range("K113").Value = 9091668
$ws.Range("M113").Value = -9089498

$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550

$ws.Range("H141").Value = 382593.94
$ws.Range("J141").Value = 479055.72
$ws.Range("L141").Value = 479055.72
$ws.Range("N141").Value = -489415.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 562.3333
$ws.Range("I46").Value = 528
$ws.Range("K46").Value = 1584
$ws.Range("M46").Value = -1493

$ws.Range("H109").Value = 3697.375
$ws.Range("I109").Value = 1116
$ws.Range("K109").Value = 3348
$ws.Range("M109").Value = -2308

$ws.Range("H122").Value = 17599.143
$ws.Range("I122").Value = 42225
$ws.Range("J122").Value = 2444.7693
$ws.Range("K122").Value = 380025
$ws.Range("L122").Value = 22002.9237
$ws.Range("M122").Value = -377575
$ws.Range("N122").Value = -26902.9237

$ws.Range("H131").Value = 3445.1458
$ws.Range("I131").Value = 3400
$ws.Range("J131").Value = 3454.175
$ws.Range("K131").Value = 10200
$ws.Range("L131").Value = 10362.525
$ws.Range("M131").Value = -5160
$ws.Range("N131").Value = -20442.525

$ws.Range("H137").Value = 8684.571
$ws.Range("I137").Value = 7175
$ws.Range("J137").Value = 9039.764999999999
$ws.Range("K137").Value = 21525
$ws.Range("L137").Value = 27119.295
$ws.Range("M137").Value = -16425
$ws.Range("N137").Value = -37319.295

$ws.Range("H140").Value = 3444
$ws.Range("I140").Value = 1497.579
$ws.Range("J140").Value = 8066.75
$ws.Range("K140").Value = 4492.737
$ws.Range("L140").Value = 24200.25
$ws.Range("M140").Value = 687.2629999999999
$ws.Range("N140").Value = -34560.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 1180294

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7880.528
$ws.Range("I7").Value = 6422.7085
$ws.Range("K7").Value = 6422.7085
$ws.Range("M7").Value = -6310.7085

$ws.Range("H16").Value = 2530.0356
$ws.Range("I16").Value = 2630.3684
$ws.Range("K16").Value = 2630.3684
$ws.Range("M16").Value = -2460.3684

$ws.Range("H22").Value = 13202055
$ws.Range("J22").Value = 2706.8572
$ws.Range("L22").Value = 2706.8572
$ws.Range("N22").Value = -3296.8572

$ws.Range("H27").Value = 13202055
$ws.Range("J27").Value = 2706.8572
$ws.Range("L27").Value = 2706.8572
$ws.Range("N27").Value = -2920.8572

$ws.Range("H40").Value = 6075.5884
$ws.Range("I40").Value = 4926.5454
$ws.Range("J40").Value = 8182.1665
$ws.Range("K40").Value = 4926.5454
$ws.Range("L40").Value = 8182.1665
$ws.Range("M40").Value = -4790.5454
$ws.Range("N40").Value = -8454.166499999999

$ws.Range("H46").Value = 676.8889
$ws.Range("I46").Value = 499
$ws.Range("J46").Value = 1299.5
$ws.Range("K46").Value = 499
$ws.Range("L46").Value = 1299.5
$ws.Range("M46").Value = -311
$ws.Range("N46").Value = -1675.5

$ws.Range("H61").Value = 4754.1177
$ws.Range("I61").Value = 4662.6924
$ws.Range("K61").Value = 4662.6924
$ws.Range("M61").Value = -4460.6924

$ws.Range("H113").Value = 4754.1177
$ws.Range("I113").Value = 4662.6924
$ws.Range("K113").Value = 4662.6924
$ws.Range("M113").Value = -2492.6924

$ws.Range("H126").Value = 7880.528
$ws.Range("I126").Value = 6422.7085
$ws.Range("K126").Value = 19268.1255
$ws.Range("M126").Value = -16798.1255

$ws.Range("H136").Value = 4100.0386
$ws.Range("I136").Value = 3290
$ws.Range("K136").Value = 9870
$ws.Range("M136").Value = -7320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 36616.5
$ws.Range("J15").Value = 36616.5
$ws.Range("L15").Value = 36616.5
$ws.Range("N15").Value = -37192.5

$ws.Range("H31").Value = 12000
$ws.Range("I31").Value = 12000
$ws.Range("K31").Value = 12000
$ws.Range("M31").Value = -11652

$ws.Range("H122").Value = 3533.4
$ws.Range("I122").Value = 3462.4614
$ws.Range("K122").Value = 10387.3842
$ws.Range("M122").Value = -7937.3842

$ws.Range("H132").Value = 258093.67
$ws.Range("I132").Value = 1713.4117
$ws.Range("J132").Value = 2001479.4
$ws.Range("K132").Value = 5140.2351
$ws.Range("L132").Value = 6004438.199999999
$ws.Range("M132").Value = -2610.2351
$ws.Range("N132").Value = -6009498.199999999
